# Weekly update: insert a new daily price record for "Papa" (Patagonia,
# 1a (guarda)) at row 353 of the "Feria Lagunitas de Puerto Montt" sheet.
# Inserting (rather than appending) shifts every existing record from the
# old row 353 down by one row, which matches the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 353; everything below shifts down by one.
$ws.Rows.Item(353).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A353").Value = 4
$ws.Range("B353").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C353").Value = "Los Lagos"
$ws.Range("D353").Value = 44782
$ws.Range("E353").Value = 10
$ws.Range("F353").Value = 100114001
$ws.Range("G353").Value = "Papa"
$ws.Range("H353").Value = "Patagonia"
$ws.Range("I353").Value = "1a (guarda)"
$ws.Range("J353").Value = 600
$ws.Range("K353").Value = 8000
$ws.Range("L353").Value = 8000
$ws.Range("M353").Value = 8000
$ws.Range("N353").Value = '$/saco 25 kilos'
$ws.Range("O353").Value = "Provincia de Llanquihue"
$ws.Range("P353").Value = 320
$ws.Range("Q353").Value = 25
$ws.Range("R353").Value = "Hortaliza"
